$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "26÷2="
$t.Cell(1,2).Range.Text  = "54÷6="
$t.Cell(1,3).Range.Text  = "70÷9="
$t.Cell(1,4).Range.Text  = "37÷2="
$t.Cell(1,5).Range.Text  = "87÷7="

$t.Cell(5,1).Range.Text  = "93÷7="
$t.Cell(5,2).Range.Text  = "71÷3="
$t.Cell(5,3).Range.Text  = "28÷6="
$t.Cell(5,4).Range.Text  = "42÷2="
$t.Cell(5,5).Range.Text  = "35÷6="

$t.Cell(9,1).Range.Text  = "61÷3="
$t.Cell(9,2).Range.Text  = "27÷7="
$t.Cell(9,3).Range.Text  = "53÷2="
$t.Cell(9,4).Range.Text  = "98÷3="
$t.Cell(9,5).Range.Text  = "97÷8="

$t.Cell(13,1).Range.Text = "27÷6="
$t.Cell(13,2).Range.Text = "75÷7="
$t.Cell(13,3).Range.Text = "68÷4="
$t.Cell(13,4).Range.Text = "32÷5="
$t.Cell(13,5).Range.Text = "32÷5="

$t.Cell(17,1).Range.Text = "24÷4="
$t.Cell(17,2).Range.Text = "90÷9="
$t.Cell(17,3).Range.Text = "57÷8="
$t.Cell(17,4).Range.Text = "44÷8="
$t.Cell(17,5).Range.Text = "62÷9="
